$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 2506.25
$ws.Cells.Item(98, 10).Value = 4550
$ws.Cells.Item(98, 12).Value = 4550
$ws.Cells.Item(98, 14).Value = -7546

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(111, 8).Value = 22225138
$ws.Cells.Item(111, 9).Value = 27779838
$ws.Cells.Item(111, 11).Value = 83339514
$ws.Cells.Item(111, 13).Value = -83336447

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 23812376
$ws.Cells.Item(113, 9).Value = 38463790
$ws.Cells.Item(113, 11).Value = 38463790
$ws.Cells.Item(113, 13).Value = -38460536

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 2506.25
$ws.Cells.Item(122, 10).Value = 4550
$ws.Cells.Item(122, 12).Value = 13650
$ws.Cells.Item(122, 14).Value = -18550

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(133, 8).Value = 58333.332
$ws.Cells.Item(133, 10).Value = 58333.332
$ws.Cells.Item(133, 12).Value = 58333.332
$ws.Cells.Item(133, 14).Value = -68453.33199999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(134, 8).Value = 49348.082
$ws.Cells.Item(134, 10).Value = 49348.082
$ws.Cells.Item(134, 12).Value = 49348.082
$ws.Cells.Item(134, 14).Value = -59488.082

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 287230.28
$ws.Cells.Item(135, 9).Value = 295664.12
$ws.Cells.Item(135, 10).Value = 480
$ws.Cells.Item(135, 11).Value = 2660977.08
$ws.Cells.Item(135, 12).Value = 4320
$ws.Cells.Item(135, 13).Value = -2658442.08
$ws.Cells.Item(135, 14).Value = -9390

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(136, 8).Value = 59231.285
$ws.Cells.Item(136, 10).Value = 59231.285
$ws.Cells.Item(136, 12).Value = 59231.285
$ws.Cells.Item(136, 14).Value = -69431.285

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1672.7441
$ws.Cells.Item(137, 9).Value = 1278.8518
$ws.Cells.Item(137, 10).Value = 2337.4375
$ws.Cells.Item(137, 11).Value = 3836.5554
$ws.Cells.Item(137, 12).Value = 7012.3125
$ws.Cells.Item(137, 13).Value = -1286.5554
$ws.Cells.Item(137, 14).Value = -12112.3125

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 2276.0286
$ws.Cells.Item(141, 9).Value = 2055.5715
$ws.Cells.Item(141, 10).Value = 3157.8572
$ws.Cells.Item(141, 11).Value = 6166.7145
$ws.Cells.Item(141, 12).Value = 9473.571599999999
$ws.Cells.Item(141, 13).Value = -986.7145
$ws.Cells.Item(141, 14).Value = -19833.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 180
$ws.Cells.Item(4, 9).Value = 180
$ws.Cells.Item(4, 11).Value = 180
$ws.Cells.Item(4, 13).Value = -64

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(76, 8).Value = 32444
$ws.Cells.Item(76, 10).Value = 32444
$ws.Cells.Item(76, 12).Value = 32444
$ws.Cells.Item(76, 14).Value = -33120

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(79, 8).Value = 32444
$ws.Cells.Item(79, 10).Value = 32444
$ws.Cells.Item(79, 12).Value = 32444
$ws.Cells.Item(79, 14).Value = -34784

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(134, 8).Value = 38999.7
$ws.Cells.Item(134, 10).Value = 38999.7
$ws.Cells.Item(134, 12).Value = 38999.7
$ws.Cells.Item(134, 14).Value = -49139.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(141, 8).Value = 58576.332
$ws.Cells.Item(141, 10).Value = 58576.332
$ws.Cells.Item(141, 12).Value = 58576.332
$ws.Cells.Item(141, 14).Value = -68936.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(123, 8).Value = 33745
$ws.Cells.Item(123, 10).Value = 33745
$ws.Cells.Item(123, 12).Value = 33745
$ws.Cells.Item(123, 14).Value = -43545

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(18, 8).Value = 34680
$ws.Cells.Item(18, 10).Value = 34680
$ws.Cells.Item(18, 12).Value = 34680
$ws.Cells.Item(18, 14).Value = -35140

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 1645.75
$ws.Cells.Item(99, 9).Value = 1595.4286
$ws.Cells.Item(99, 10).Value = 1998
$ws.Cells.Item(99, 11).Value = 1595.4286
$ws.Cells.Item(99, 12).Value = 1998
$ws.Cells.Item(99, 13).Value = -97.42859999999996
$ws.Cells.Item(99, 14).Value = -4994

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 2043.9412
$ws.Cells.Item(122, 9).Value = 1219.4166
$ws.Cells.Item(122, 11).Value = 3658.2498
$ws.Cells.Item(122, 13).Value = -1208.2498

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 1645.75
$ws.Cells.Item(126, 9).Value = 1595.4286
$ws.Cells.Item(126, 10).Value = 1998
$ws.Cells.Item(126, 11).Value = 4786.2858
$ws.Cells.Item(126, 12).Value = 5994
$ws.Cells.Item(126, 13).Value = -2316.2858
$ws.Cells.Item(126, 14).Value = -10934

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(140, 8).Value = 56661.59
$ws.Cells.Item(140, 10).Value = 56661.59
$ws.Cells.Item(140, 12).Value = 56661.59
$ws.Cells.Item(140, 14).Value = -67021.59

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 878.29
$ws.Cells.Item(131, 9).Value = 480.8889
$ws.Cells.Item(131, 10).Value = 917.5934
$ws.Cells.Item(131, 11).Value = 1442.6667
$ws.Cells.Item(131, 12).Value = 2752.7802
$ws.Cells.Item(131, 13).Value = 3597.3333
$ws.Cells.Item(131, 14).Value = -12832.7802

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 879531
$ws.Cells.Item(132, 9).Value = 1462935.4
$ws.Cells.Item(132, 11).Value = 13166418.6
$ws.Cells.Item(132, 13).Value = -13163888.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(134, 8).Value = 2016.5927
$ws.Cells.Item(134, 9).Value = 1488.5454
$ws.Cells.Item(134, 10).Value = 4340
$ws.Cells.Item(134, 11).Value = 4465.6362
$ws.Cells.Item(134, 12).Value = 13020
$ws.Cells.Item(134, 13).Value = 604.3638000000001
$ws.Cells.Item(134, 14).Value = -23160

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 26733.824
$ws.Cells.Item(139, 9).Value = 1630.2703
$ws.Cells.Item(139, 10).Value = 336344.34
$ws.Cells.Item(139, 11).Value = 4890.810899999999
$ws.Cells.Item(139, 12).Value = 1009033.02
$ws.Cells.Item(139, 13).Value = 249.1891000000005
$ws.Cells.Item(139, 14).Value = -1019313.02

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 65282.562
$ws.Cells.Item(140, 9).Value = 112668.89
$ws.Cells.Item(140, 10).Value = 4357.2856
$ws.Cells.Item(140, 11).Value = 338006.67
$ws.Cells.Item(140, 12).Value = 13071.8568
$ws.Cells.Item(140, 13).Value = -332826.67
$ws.Cells.Item(140, 14).Value = -23431.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2448.7896
$ws.Cells.Item(102, 9).Value = 1691.4166
$ws.Cells.Item(102, 11).Value = 1691.4166
$ws.Cells.Item(102, 13).Value = -69.41660000000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1731.4814
$ws.Cells.Item(122, 9).Value = 1697.9546
$ws.Cells.Item(122, 11).Value = 5093.8638
$ws.Cells.Item(122, 13).Value = -2643.8638

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 5613.4443
$ws.Cells.Item(126, 9).Value = 3038.9
$ws.Cells.Item(126, 11).Value = 9116.700000000001
$ws.Cells.Item(126, 13).Value = -6646.700000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 7694116.5
$ws.Cells.Item(7, 9).Value = 8334876
$ws.Cells.Item(7, 10).Value = 5000
$ws.Cells.Item(7, 11).Value = 8334876
$ws.Cells.Item(7, 12).Value = 5000
$ws.Cells.Item(7, 13).Value = -8334764
$ws.Cells.Item(7, 14).Value = -5224

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 1581
$ws.Cells.Item(40, 9).Value = 1523.5625
$ws.Cells.Item(40, 10).Value = 2500
$ws.Cells.Item(40, 11).Value = 1523.5625
$ws.Cells.Item(40, 12).Value = 2500
$ws.Cells.Item(40, 13).Value = -1387.5625
$ws.Cells.Item(40, 14).Value = -2772

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1816.25
$ws.Cells.Item(46, 9).Value = 1310.5555
$ws.Cells.Item(46, 10).Value = 3333.3333
$ws.Cells.Item(46, 11).Value = 1310.5555
$ws.Cells.Item(46, 12).Value = 3333.3333
$ws.Cells.Item(46, 13).Value = -1122.5555
$ws.Cells.Item(46, 14).Value = -3709.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 2673.45
$ws.Cells.Item(122, 9).Value = 2266.6667
$ws.Cells.Item(122, 10).Value = 2847.7856
$ws.Cells.Item(122, 11).Value = 6800.000100000001
$ws.Cells.Item(122, 12).Value = 8543.356800000001
$ws.Cells.Item(122, 13).Value = -4350.000100000001
$ws.Cells.Item(122, 14).Value = -13443.3568

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 7694116.5
$ws.Cells.Item(126, 9).Value = 8334876
$ws.Cells.Item(126, 10).Value = 5000
$ws.Cells.Item(126, 11).Value = 25004628
$ws.Cells.Item(126, 12).Value = 15000
$ws.Cells.Item(126, 13).Value = -25002158
$ws.Cells.Item(126, 14).Value = -19940

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 626215.5600000001
$ws.Cells.Item(100, 9).Value = 1401.1538
$ws.Cells.Item(100, 10).Value = 3333744.8
$ws.Cells.Item(100, 11).Value = 2802.3076
$ws.Cells.Item(100, 12).Value = 6667489.6
$ws.Cells.Item(100, 13).Value = -2261.3076
$ws.Cells.Item(100, 14).Value = -6668571.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(108, 8).Value = 19800
$ws.Cells.Item(108, 10).Value = 19800
$ws.Cells.Item(108, 12).Value = 19800
$ws.Cells.Item(108, 14).Value = -27480

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 8973.379000000001
$ws.Cells.Item(122, 9).Value = 11352.762
$ws.Cells.Item(122, 10).Value = 2727.5
$ws.Cells.Item(122, 11).Value = 34058.286
$ws.Cells.Item(122, 12).Value = 8182.5
$ws.Cells.Item(122, 13).Value = -31608.286
$ws.Cells.Item(122, 14).Value = -13082.5
